# Add a new paragraph "Prueba 03" after the existing "PRUEBA 02" paragraph,
# with a collapsed "_GoBack" bookmark sitting right after the new text
# (mirrors Word's auto-added last-edit-position bookmark).

$d = $word.ActiveDocument

# New empty paragraph at the end of the document.
$d.Content.InsertParagraphAfter()

# Type the paragraph text, followed by a throwaway sentinel character so the
# insertion point used for the bookmark sits mid-paragraph (not exactly on
# the paragraph mark boundary).
$r = $d.Content
$r.Collapse(0)
$r.InsertBefore("Prueba 03X")

# Anchor a collapsed bookmark immediately before the sentinel character.
$bmPos = $r.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the sentinel character, leaving the bookmark collapsed right after
# "Prueba 03".
$sentinel = $d.Range($bmPos, $bmPos + 1)
$sentinel.Delete()
